$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 162, shifting the rest of the
# weekly data (previous rows 162-168) down to rows 164-170.
$ws.Rows.Item(162).Insert()
$ws.Rows.Item(162).Insert()

# --- New row 162: Navel Late / Primera -----------------------------------
$ws.Cells.Item(162,1).Value  = 11
$ws.Cells.Item(162,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(162,3).Value  = "Bíobío"
$ws.Cells.Item(162,4).Value  = 44509
$ws.Cells.Item(162,5).Value  = 8
$ws.Cells.Item(162,6).Value  = "Fruta"
$ws.Cells.Item(162,7).Value  = 100102
$ws.Cells.Item(162,8).Value  = "Cítricos"
$ws.Cells.Item(162,9).Value  = 100102005
$ws.Cells.Item(162,10).Value = "Naranja"
$ws.Cells.Item(162,11).Value = "Navel Late"
$ws.Cells.Item(162,12).Value = "Primera"
$ws.Cells.Item(162,13).Value = 650
$ws.Cells.Item(162,14).Value = 8000
$ws.Cells.Item(162,15).Value = 8500
$ws.Cells.Item(162,16).Value = 8269
$ws.Cells.Item(162,17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(162,18).Value = "Región de O'Higgins"
$ws.Cells.Item(162,19).Value = 551
$ws.Cells.Item(162,20).Value = 15

# --- New row 163: Olinda Valencia / Primera -------------------------------
$ws.Cells.Item(163,1).Value  = 11
$ws.Cells.Item(163,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(163,3).Value  = "Bíobío"
$ws.Cells.Item(163,4).Value  = 44509
$ws.Cells.Item(163,5).Value  = 8
$ws.Cells.Item(163,6).Value  = "Fruta"
$ws.Cells.Item(163,7).Value  = 100102
$ws.Cells.Item(163,8).Value  = "Cítricos"
$ws.Cells.Item(163,9).Value  = 100102005
$ws.Cells.Item(163,10).Value = "Naranja"
$ws.Cells.Item(163,11).Value = "Olinda Valencia"
$ws.Cells.Item(163,12).Value = "Primera"
$ws.Cells.Item(163,13).Value = 450
$ws.Cells.Item(163,14).Value = 7500
$ws.Cells.Item(163,15).Value = 8000
$ws.Cells.Item(163,16).Value = 7722
$ws.Cells.Item(163,17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(163,18).Value = "Región de O'Higgins"
$ws.Cells.Item(163,19).Value = 515
$ws.Cells.Item(163,20).Value = 15
